$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 18: Nivel/nivel -> ...17/x17
$ws.Range("A18").Value = "...17"
$ws.Range("B18").Value = "x17"

# Row 56: Menores/menores -> Menores en el hogar/menores_en_el_hogar
$ws.Range("A56").Value = "Menores en el hogar"
$ws.Range("B56").Value = "menores_en_el_hogar"
